$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B51 should become a real number (4), not a string
$ws.Range("B51").Value = 4

# Add new row 52
$ws.Range("A52").Value = "Sunsi Wu"
$ws.Range("B52").Value = "'3"
$ws.Range("B52").Style = "Normal"
$ws.Range("C52").Value = "无"
$ws.Range("D52").Value = "ACK"
$ws.Range("E52").Value = "EXP"
$ws.Range("F52").Value = "6dbc86e6-aac5-4bea-af0c-fc9177dfd16b"
$ws.Range("G52").Value = "BkJ3ibb0-_annotated.xlsx"
$ws.Range("H52").Value = "Furthermore, we have not optimized the running time of our algorithm, as it was not the focus of this work."
